$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values must stay plain text (as in the source data) even though many
# look numeric, so we force Text format before the write and clear the helper
# format afterwards to avoid leaving a residual cell style behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.936.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.255.61"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.43"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.524"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.69"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.606.05"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.34"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.248.70"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.819.96"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.31"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.17%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.09"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.61"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.41"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.95"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.113"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.960.14"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.47"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.85"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.74%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.11"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.94%  "
